$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old trailing rows (29-35) which no longer exist in the corrected data
$ws.Rows("29:35").Delete()

$ws.Range("A5").Value2 = "homens"
$ws.Range("B5").Value2 = 2.34
$ws.Range("C5").Value2 = 2.46
$ws.Range("D5").Value2 = 2.7
$ws.Range("E5").Value2 = 2.52
$ws.Range("F5").Value2 = 1.64

$ws.Range("A6").Value2 = "mulheres"
$ws.Range("B6").Value2 = 2.32
$ws.Range("C6").Value2 = 2.06
$ws.Range("D6").Value2 = 2.59
$ws.Range("E6").Value2 = 5.06
$ws.Range("F6").Value2 = 2.17

$ws.Range("A7").Value2 = "branca"
$ws.Range("B7").Value2 = 2.44
$ws.Range("C7").Value2 = 2.49
$ws.Range("D7").Value2 = 2.67
$ws.Range("E7").Value2 = 3.41
$ws.Range("F7").Value2 = 2.02

$ws.Range("A8").Value2 = "preta ou parda"
$ws.Range("B8").Value2 = 2.48
$ws.Range("C8").Value2 = 2.18
$ws.Range("D8").Value2 = 2.77
$ws.Range("E8").Value2 = 3.53
$ws.Range("F8").Value2 = 1.87

$ws.Range("A9").Value2 = "preta"
$ws.Range("B9").Value2 = 5.63
$ws.Range("C9").Value2 = 5.58
$ws.Range("D9").Value2 = 6.55
$ws.Range("E9").Value2 = 8.06
$ws.Range("F9").Value2 = 4.49

$ws.Range("A10").Value2 = "parda"
$ws.Range("B10").Value2 = 2.65
$ws.Range("C10").Value2 = 2.35
$ws.Range("D10").Value2 = 3.01
$ws.Range("E10").Value2 = 3.57
$ws.Range("F10").Value2 = 2.06

$ws.Range("A11").Value2 = "10 a 24 anos"
$ws.Range("B11").Value2 = 4.09
$ws.Range("C11").Value2 = 2.32
$ws.Range("D11").Value2 = 4.15
$ws.Range("E11").Value2 = 4.59
$ws.Range("F11").Value2 = 4.06

$ws.Range("A12").Value2 = "25 a 34 anos"
$ws.Range("B12").Value2 = 3.3
$ws.Range("C12").Value2 = 3.24
$ws.Range("D12").Value2 = 3.71
$ws.Range("E12").Value2 = 4.63
$ws.Range("F12").Value2 = 2.74

$ws.Range("A13").Value2 = "35 a 39 anos"
$ws.Range("B13").Value2 = 5
$ws.Range("C13").Value2 = 5.28
$ws.Range("D13").Value2 = 5.2
$ws.Range("E13").Value2 = 5.77
$ws.Range("F13").Value2 = 3.67

$ws.Range("A14").Value2 = "40 a 49 anos"
$ws.Range("B14").Value2 = 3.48
$ws.Range("C14").Value2 = 4.33
$ws.Range("D14").Value2 = 4.27
$ws.Range("E14").Value2 = 4.79
$ws.Range("F14").Value2 = 2.64

$ws.Range("A15").Value2 = "50 anos ou mais"
$ws.Range("B15").Value2 = 2.72
$ws.Range("C15").Value2 = 4.29
$ws.Range("D15").Value2 = 3.59
$ws.Range("E15").Value2 = 5.22
$ws.Range("F15").Value2 = 2.15

$ws.Range("A16").Value2 = "sem instrução"
$ws.Range("B16").Value2 = 5.23
$ws.Range("C16").Value2 = 9.109999999999999
$ws.Range("D16").Value2 = 6.96
$ws.Range("E16").Value2 = 9.74
$ws.Range("F16").Value2 = 3.59

$ws.Range("A17").Value2 = "fundamental incompleto ou equivalente"
$ws.Range("B17").Value2 = 3.01
$ws.Range("C17").Value2 = 3.1
$ws.Range("D17").Value2 = 3.6
$ws.Range("E17").Value2 = 4.02
$ws.Range("F17").Value2 = 2.17

$ws.Range("A18").Value2 = "fundamental completo ou equivalente"
$ws.Range("B18").Value2 = 5.41
$ws.Range("C18").Value2 = 4.85
$ws.Range("D18").Value2 = 6.05
$ws.Range("E18").Value2 = 6.07
$ws.Range("F18").Value2 = 4.05

$ws.Range("A19").Value2 = "médio incompleto ou equivalente"
$ws.Range("B19").Value2 = 6.35
$ws.Range("C19").Value2 = 4.22
$ws.Range("D19").Value2 = 6.5
$ws.Range("E19").Value2 = 8.26
$ws.Range("F19").Value2 = 5.75

$ws.Range("A20").Value2 = "médio completo ou equivalente"
$ws.Range("B20").Value2 = 3.05
$ws.Range("C20").Value2 = 3.01
$ws.Range("D20").Value2 = 3.61
$ws.Range("E20").Value2 = 4.34
$ws.Range("F20").Value2 = 2.81

$ws.Range("A21").Value2 = "superior incompleto ou equivalente"
$ws.Range("B21").Value2 = 5.8
$ws.Range("C21").Value2 = 4.49
$ws.Range("D21").Value2 = 6.25
$ws.Range("E21").Value2 = 9.630000000000001
$ws.Range("F21").Value2 = 5.65

$ws.Range("A22").Value2 = "superior completo ou equivalente"
$ws.Range("B22").Value2 = 4.2
$ws.Range("C22").Value2 = 4.81
$ws.Range("D22").Value2 = 4.48
$ws.Range("E22").Value2 = 7.89
$ws.Range("F22").Value2 = 3.9

$ws.Range("A23").Value2 = "total (3)(4)"
$ws.Range("B23").Value2 = 1.8
$ws.Range("C23").Value2 = 1.69
$ws.Range("D23").Value2 = 2
$ws.Range("E23").Value2 = 2.63
$ws.Range("F23").Value2 = 1.45

$ws.Range("A24").Value2 = "sem rendimento a menos de 1/4 do salário mínimo (3) (5)"
$ws.Range("B24").Value2 = 5.95
$ws.Range("C24").Value2 = 6.61
$ws.Range("D24").Value2 = 7.72
$ws.Range("E24").Value2 = 9.73
$ws.Range("F24").Value2 = 4.98

$ws.Range("A25").Value2 = "1/4 a menos de 1/2 salário mínimo (3)"
$ws.Range("B25").Value2 = 4.85
$ws.Range("C25").Value2 = 4.54
$ws.Range("D25").Value2 = 5.83
$ws.Range("E25").Value2 = 5.59
$ws.Range("F25").Value2 = 3.35

$ws.Range("A26").Value2 = "1/2 a menos de 1 salário mínimo (3)"
$ws.Range("B26").Value2 = 3.5
$ws.Range("C26").Value2 = 3.22
$ws.Range("D26").Value2 = 4.28
$ws.Range("E26").Value2 = 4.25
$ws.Range("F26").Value2 = 2.91

$ws.Range("A27").Value2 = "1 a menos de 2 salários mínimos (3) "
$ws.Range("B27").Value2 = 3.1
$ws.Range("C27").Value2 = 3.06
$ws.Range("D27").Value2 = 3.49
$ws.Range("E27").Value2 = 4.15
$ws.Range("F27").Value2 = 2.47

$ws.Range("A28").Value2 = "2 salários mínimos ou mais (3) "
$ws.Range("B28").Value2 = 3.03
$ws.Range("C28").Value2 = 3.32
$ws.Range("D28").Value2 = 3.31
$ws.Range("E28").Value2 = 5.15
$ws.Range("F28").Value2 = 2.76
